$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Authors column (E) with refreshed author/affiliation-count strings
$ws.Cells.Item(5, 5).Value = '[Chaolin%Huang%NULL%0,  Yeming%Wang%NULL%4,  Xingwang%Li%NULL%4,  Lili%Ren%NULL%2,  Jianping%Zhao%NULL%2,  Yi%Hu%NULL%4,  Li%Zhang%NULL%4,  Guohui%Fan%NULL%4,  Jiuyang%Xu%NULL%4,  Xiaoying%Gu%NULL%4,  Zhenshun%Cheng%NULL%2,  Ting%Yu%NULL%8,  Jiaan%Xia%NULL%2,  Yuan%Wei%NULL%7,  Wenjuan%Wu%NULL%2,  Xuelei%Xie%NULL%2,  Wen%Yin%NULL%3,  Hui%Li%NULL%5,  Min%Liu%NULL%2,  Yan%Xiao%NULL%3,  Hong%Gao%NULL%3,  Li%Guo%NULL%3,  Jungang%Xie%NULL%2,  Guangfa%Wang%NULL%2,  Rongmeng%Jiang%NULL%2,  Zhancheng%Gao%NULL%3,  Qi%Jin%NULL%3,  Jianwei%Wang%wangjw28@163.com%2,  Bin%Cao%caobin_ben@163.com%5]'
$ws.Cells.Item(6, 5).Value = '[Yan%Deng%NULL%5,  Wei%Liu%NULL%2,  Kui%Liu%NULL%2,  Yuan-Yuan%Fang%NULL%2,  Jin%Shang%NULL%1,  Ling%Zhou%NULL%1,  Ke%Wang%NULL%1,  Fan%Leng%NULL%1,  Shuang%Wei%NULL%1,  Lei%Chen%NULL%1,  Hui-Guo%Liu%NULL%2,  Pei-Fang%Wei%NULL%6,  Pei-Fang%Wei%NULL%0]'
$ws.Cells.Item(7, 5).Value = '[Yichun%Cheng%NULL%1,  Ran%Luo%NULL%1,  Kun%Wang%NULL%2,  Meng%Zhang%NULL%1,  Zhixiang%Wang%NULL%1,  Lei%Dong%NULL%1,  Junhua%Li%NULL%2,  Ying%Yao%NULL%1,  Shuwang%Ge%geshuwang@tjh.tjmu.edu.cn%1,  Gang%Xu%xugang@tjh.tjmu.edu.cn%1]'
$ws.Cells.Item(8, 5).Value = '[Carly%Eastin%NULL%2,  Travis%Eastin%NULL%1]'
$ws.Cells.Item(11, 5).Value = '[Luwen%Wang%NULL%2,  Xun%Li%NULL%2,  Hui%Chen%NULL%1,  Shaonan%Yan%NULL%1,  Dong%Li%NULL%1,  Yan%Li%NULL%1,  Zuojiong%Gong%NULL%1]'
$ws.Cells.Item(12, 5).Value = '[Xiao-Wei%Xu%NULL%0,  Xiao-Xin%Wu%NULL%3,  Xian-Gao%Jiang%NULL%3,  Kai-Jin%Xu%NULL%3,  Ling-Jun%Ying%NULL%3,  Chun-Lian%Ma%NULL%3,  Shi-Bo%Li%NULL%3,  Hua-Ying%Wang%NULL%3,  Sheng%Zhang%NULL%3,  Hai-Nv%Gao%NULL%3,  Ji-Fang%Sheng%NULL%3,  Hong-Liu%Cai%NULL%3,  Yun-Qing%Qiu%NULL%3,  Lan-Juan%Li%NULL%3]'
$ws.Cells.Item(14, 5).Value = '[Gemin%Zhang%NULL%3,  Jie%Zhang%945128911@qq.com%1,  Bowen%Wang%NULL%1,  Xionglin%Zhu%NULL%1,  Qiang%Wang%NULL%1,  Shiming%Qiu%NULL%1]'
$ws.Cells.Item(15, 5).Value = '[Xiaoli%Zhang%NULL%2,  Huan%Cai%NULL%1,  Jianhua%Hu%NULL%1,  Jiangshan%Lian%NULL%1,  Jueqing%Gu%NULL%1,  Shanyan%Zhang%NULL%1,  Chanyuan%Ye%NULL%1,  Yingfeng%Lu%NULL%1,  Ciliang%Jin%NULL%1,  Guodong%Yu%NULL%1,  Hongyu%Jia%NULL%1,  Yimin%Zhang%NULL%1,  Jifang%Sheng%jifang_sheng@zju.edu.cn%1,  Lanjuan%Li%ljli@zju.edu.cn%2,  Yida%Yang%yidayang65@zju.edu.cn%1]'
$ws.Cells.Item(16, 5).Value = '[Fei%Zhou%NULL%0,  Ting%Yu%NULL%0,  Ronghui%Du%NULL%2,  Guohui%Fan%NULL%0,  Ying%Liu%NULL%5,  Zhibo%Liu%NULL%2,  Jie%Xiang%NULL%4,  Yeming%Wang%NULL%0,  Bin%Song%NULL%2,  Xiaoying%Gu%NULL%0,  Lulu%Guan%NULL%2,  Yuan%Wei%NULL%0,  Hui%Li%NULL%0,  Xudong%Wu%NULL%2,  Jiuyang%Xu%NULL%0,  Shengjin%Tu%NULL%2,  Yi%Zhang%NULL%2,  Hua%Chen%NULL%2,  Bin%Cao%NULL%0]'

# Add new column I ("Other found locations") values
$ws.Cells.Item(1, 9).Value = 'Other found locations'
$ws.Cells.Item(2, 9).Value = ''
$ws.Cells.Item(3, 9).Value = ''
$ws.Cells.Item(4, 9).Value = ''
$ws.Cells.Item(5, 9).Value = '_PMC_elsevier'
$ws.Cells.Item(6, 9).Value = '_PMC'
$ws.Cells.Item(7, 9).Value = '_PMC_elsevier'
$ws.Cells.Item(8, 9).Value = '_PMC_elsevier'
$ws.Cells.Item(9, 9).Value = ''
$ws.Cells.Item(10, 9).Value = '_PMC'
$ws.Cells.Item(11, 9).Value = '_PMC'
$ws.Cells.Item(12, 9).Value = '_PMC'
$ws.Cells.Item(13, 9).Value = ''
$ws.Cells.Item(14, 9).Value = '_PMC_Springer'
$ws.Cells.Item(15, 9).Value = '_PMC_elsevier'
$ws.Cells.Item(16, 9).Value = '_PMC_elsevier'
$ws.Cells.Item(17, 9).Value = ''
